# Update TPM-derived values in the LR-pairs sheet (Cxcl16-Cxcr6)
# per "update scripts wuth new tpm" commit: refresh G/H/I/J (ligand),
# M/N/O/P (receptor) and Q/R/S/T (edge) derived columns for rows 2-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.661288333333333
$ws.Range("H2").Value = 4.983865
$ws.Range("I2").Value = 0.1812815212055838
$ws.Range("J2").Value = 0.1812815212055839
$ws.Range("M2").Value = 1.046760666666667
$ws.Range("N2").Value = 3.140282
$ws.Range("O2").Value = 0.2608060558711016
$ws.Range("P2").Value = 0.2608060558711016
$ws.Range("Q2").Value = 1.738971283325556
$ws.Range("R2").Value = 15.65074154993
$ws.Range("S2").Value = 0.04727931854794179
$ws.Range("T2").Value = 0.0472793185479418

# Row 3
$ws.Range("G3").Value = 1.661288333333333
$ws.Range("H3").Value = 4.983865
$ws.Range("I3").Value = 0.1812815212055838
$ws.Range("J3").Value = 0.1812815212055839
$ws.Range("O3").Value = 0.4138468436871376
$ws.Range("P3").Value = 0.4138468436871375
$ws.Range("Q3").Value = 2.759398260378334
$ws.Range("R3").Value = 24.834584343405
$ws.Range("S3").Value = 0.07502278536973377
$ws.Range("T3").Value = 0.07502278536973377

# Row 4
$ws.Range("G4").Value = 1.661288333333333
$ws.Range("H4").Value = 4.983865
$ws.Range("I4").Value = 0.1812815212055838
$ws.Range("J4").Value = 0.1812815212055839
$ws.Range("O4").Value = 0.3253471004417607
$ws.Range("P4").Value = 0.3253471004417607
$ws.Range("Q4").Value = 2.169310305666666
$ws.Range("R4").Value = 19.523792751
$ws.Range("S4").Value = 0.05897941728790826
$ws.Range("T4").Value = 0.05897941728790827

# Row 5
$ws.Range("I5").Value = 0.7639341119847463
$ws.Range("J5").Value = 0.7639341119847464
$ws.Range("M5").Value = 1.046760666666667
$ws.Range("N5").Value = 3.140282
$ws.Range("O5").Value = 0.2608060558711016
$ws.Range("P5").Value = 0.2608060558711016
$ws.Range("Q5").Value = 7.328157190316889
$ws.Range("R5").Value = 65.95341471285201
$ws.Range("S5").Value = 0.1992386426921341
$ws.Range("T5").Value = 0.1992386426921341

# Row 6
$ws.Range("I6").Value = 0.7639341119847463
$ws.Range("J6").Value = 0.7639341119847464
$ws.Range("O6").Value = 0.4138468436871376
$ws.Range("P6").Value = 0.4138468436871375
$ws.Range("S6").Value = 0.3161517210298235
$ws.Range("T6").Value = 0.3161517210298235

# Row 7
$ws.Range("I7").Value = 0.7639341119847463
$ws.Range("J7").Value = 0.7639341119847464
$ws.Range("O7").Value = 0.3253471004417607
$ws.Range("P7").Value = 0.3253471004417607
$ws.Range("S7").Value = 0.2485437482627885
$ws.Range("T7").Value = 0.2485437482627886

# Row 8
$ws.Range("I8").Value = 0.05478436680966978
$ws.Range("J8").Value = 0.0547843668096698
$ws.Range("M8").Value = 1.046760666666667
$ws.Range("N8").Value = 3.140282
$ws.Range("O8").Value = 0.2608060558711016
$ws.Range("P8").Value = 0.2608060558711016
$ws.Range("Q8").Value = 0.5255275883808889
$ws.Range("R8").Value = 4.729748295428
$ws.Range("S8").Value = 0.01428809463102566
$ws.Range("T8").Value = 0.01428809463102567

# Row 9
$ws.Range("I9").Value = 0.05478436680966978
$ws.Range("J9").Value = 0.0547843668096698
$ws.Range("O9").Value = 0.4138468436871376
$ws.Range("P9").Value = 0.4138468436871375
$ws.Range("S9").Value = 0.02267233728758022
$ws.Range("T9").Value = 0.02267233728758022

# Row 10
$ws.Range("I10").Value = 0.05478436680966978
$ws.Range("J10").Value = 0.0547843668096698
$ws.Range("O10").Value = 0.3253471004417607
$ws.Range("P10").Value = 0.3253471004417607
$ws.Range("S10").Value = 0.0178239348910639
$ws.Range("T10").Value = 0.0178239348910639
